$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HUSA")
$ws.Columns("D:D").Insert()
